$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2:C118").Value = 45206
